$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item(1).Name = "hojaejemplo11"
$wb.Worksheets.Item(2).Name = "hojaejemplo12"
$wb.Worksheets.Item(3).Name = "hojaejemplo13"

# Mirror column A into column B on each sheet
$lastRows = @(4, 5, 2)
for ($s = 1; $s -le 3; $s++) {
    $ws = $wb.Worksheets.Item($s)
    $lastRow = $lastRows[$s - 1]
    for ($r = 1; $r -le $lastRow; $r++) {
        $val = $ws.Cells.Item($r, 1).Value2
        $ws.Cells.Item($r, 2).Value2 = $val
    }
    $ws.Range("B1:B2").Select()
}
